$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-5) are cyclically rotated for columns D, L, M, N, O, P, S:
#   new Row2 = old Row5, new Row3 = old Row2, new Row4 = old Row3, new Row5 = old Row4

$ws.Range("D2").Value = 44991
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 6000
$ws.Range("P2").Value = 6000
$ws.Range("S2").Value = 3000

$ws.Range("D3").Value = 44995
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 5500
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 5750
$ws.Range("S3").Value = 2875

$ws.Range("D4").Value = 45008
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("S4").Value = 3500

$ws.Range("D5").Value = 45008
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 6000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("S5").Value = 3000
